$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 9.108069666666667
$ws.Range("H2").Value = 27.324209
$ws.Range("I2").Value = 0.00155006418458712
$ws.Range("J2").Value = 0.00155006418458712
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 98.946724
$ws.Range("N2").Value = 296.840172
$ws.Range("O2").Value = 0.2098009692989996
$ws.Range("P2").Value = 0.2098009692989996
$ws.Range("Q2").Value = 901.2136554804388
$ws.Range("R2").Value = 8110.922899323948
$ws.Range("S2").Value = 0.0003252049684020411
$ws.Range("T2").Value = 0.0003252049684020411

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 9.108069666666667
$ws.Range("H3").Value = 27.324209
$ws.Range("I3").Value = 0.00155006418458712
$ws.Range("J3").Value = 0.00155006418458712
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 163.0062356666667
$ws.Range("N3").Value = 489.018707
$ws.Range("O3").Value = 0.345629090707923
$ws.Range("P3").Value = 0.3456290907079231
$ws.Range("Q3").Value = 1484.672150553085
$ws.Range("R3").Value = 13362.04935497776
$ws.Range("S3").Value = 0.0005357472746577642
$ws.Range("T3").Value = 0.0005357472746577644

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 9.108069666666667
$ws.Range("H4").Value = 27.324209
$ws.Range("I4").Value = 0.00155006418458712
$ws.Range("J4").Value = 0.00155006418458712
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 65.39610666666668
$ws.Range("N4").Value = 196.18832
$ws.Range("O4").Value = 0.1386621609326595
$ws.Range("P4").Value = 0.1386621609326595
$ws.Range("Q4").Value = 595.6322954487646
$ws.Range("R4").Value = 5360.690659038881
$ws.Range("S4").Value = 0.0002149352494191707
$ws.Range("T4").Value = 0.0002149352494191708

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 9.108069666666667
$ws.Range("H5").Value = 27.324209
$ws.Range("I5").Value = 0.00155006418458712
$ws.Range("J5").Value = 0.00155006418458712
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 144.2727966666667
$ws.Range("N5").Value = 432.81839
$ws.Range("O5").Value = 0.3059077790604178
$ws.Range("P5").Value = 0.3059077790604179
$ws.Range("Q5").Value = 1314.046683044834
$ws.Range("R5").Value = 11826.42014740351
$ws.Range("S5").Value = 0.0004741766921081433
$ws.Range("T5").Value = 0.0004741766921081433

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 5771.873535333333
$ws.Range("H6").Value = 17315.620606
$ws.Range("I6").Value = 0.9822909543423312
$ws.Range("J6").Value = 0.9822909543423313
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 98.946724
$ws.Range("N6").Value = 296.840172
$ws.Range("O6").Value = 0.2098009692989996
$ws.Range("P6").Value = 0.2098009692989996
$ws.Range("Q6").Value = 571107.9776635316
$ws.Range("R6").Value = 5139971.798971784
$ws.Range("S6").Value = 0.2060855943546604
$ws.Range("T6").Value = 0.2060855943546605

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 5771.873535333333
$ws.Range("H7").Value = 17315.620606
$ws.Range("I7").Value = 0.9822909543423312
$ws.Range("J7").Value = 0.9822909543423313
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 163.0062356666667
$ws.Range("N7").Value = 489.018707
$ws.Range("O7").Value = 0.345629090707923
$ws.Range("P7").Value = 0.3456290907079231
$ws.Range("Q7").Value = 940851.3777387418
$ws.Range("R7").Value = 8467662.399648678
$ws.Range("S7").Value = 0.3395083293599578
$ws.Range("T7").Value = 0.3395083293599579

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 5771.873535333333
$ws.Range("H8").Value = 17315.620606
$ws.Range("I8").Value = 0.9822909543423312
$ws.Range("J8").Value = 0.9822909543423313
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 65.39610666666668
$ws.Range("N8").Value = 196.18832
$ws.Range("O8").Value = 0.1386621609326595
$ws.Range("P8").Value = 0.1386621609326595
$ws.Range("Q8").Value = 377458.0573831692
$ws.Range("R8").Value = 3397122.516448522
$ws.Range("S8").Value = 0.136206586393712
$ws.Range("T8").Value = 0.136206586393712

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 5771.873535333333
$ws.Range("H9").Value = 17315.620606
$ws.Range("I9").Value = 0.9822909543423312
$ws.Range("J9").Value = 0.9822909543423313
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 144.2727966666667
$ws.Range("N9").Value = 432.81839
$ws.Range("O9").Value = 0.3059077790604178
$ws.Range("P9").Value = 0.3059077790604179
$ws.Range("Q9").Value = 832724.3369488604
$ws.Range("R9").Value = 7494519.032539745
$ws.Range("S9").Value = 0.3004904442340008
$ws.Range("T9").Value = 0.3004904442340009

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 1.272029666666667
$ws.Range("H10").Value = 3.816089
$ws.Range("I10").Value = 0.0002164813950916887
$ws.Range("J10").Value = 0.0002164813950916887
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 98.946724
$ws.Range("N10").Value = 296.840172
$ws.Range("O10").Value = 0.2098009692989996
$ws.Range("P10").Value = 0.2098009692989996
$ws.Range("Q10").Value = 125.8631683474787
$ws.Range("R10").Value = 1132.768515127308
$ws.Range("S10").Value = 0.00004541800652543598
$ws.Range("T10").Value = 0.00004541800652543599

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 1.272029666666667
$ws.Range("H11").Value = 3.816089
$ws.Range("I11").Value = 0.0002164813950916887
$ws.Range("J11").Value = 0.0002164813950916887
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 163.0062356666667
$ws.Range("N11").Value = 489.018707
$ws.Range("O11").Value = 0.345629090707923
$ws.Range("P11").Value = 0.3456290907079231
$ws.Range("Q11").Value = 207.3487676196581
$ws.Range("R11").Value = 1866.138908576923
$ws.Range("S11").Value = 0.00007482226774072299
$ws.Range("T11").Value = 0.00007482226774072302

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 1.272029666666667
$ws.Range("H12").Value = 3.816089
$ws.Range("I12").Value = 0.0002164813950916887
$ws.Range("J12").Value = 0.0002164813950916887
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 65.39610666666668
$ws.Range("N12").Value = 196.18832
$ws.Range("O12").Value = 0.1386621609326595
$ws.Range("P12").Value = 0.1386621609326595
$ws.Range("Q12").Value = 83.18578776449779
$ws.Range("R12").Value = 748.6720898804801
$ws.Range("S12").Value = 0.00003001777804513037
$ws.Range("T12").Value = 0.00003001777804513038

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 1.272029666666667
$ws.Range("H13").Value = 3.816089
$ws.Range("I13").Value = 0.0002164813950916887
$ws.Range("J13").Value = 0.0002164813950916887
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 144.2727966666667
$ws.Range("N13").Value = 432.81839
$ws.Range("O13").Value = 0.3059077790604178
$ws.Range("P13").Value = 0.3059077790604179
$ws.Range("Q13").Value = 183.5192774529678
$ws.Range("R13").Value = 1651.67349707671
$ws.Range("S13").Value = 0.00006622334278039933
$ws.Range("T13").Value = 0.00006622334278039934

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 93.67702500000001
$ws.Range("H14").Value = 281.031075
$ws.Range("I14").Value = 0.01594250007799006
$ws.Range("J14").Value = 0.01594250007799006
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 98.946724
$ws.Range("N14").Value = 296.840172
$ws.Range("O14").Value = 0.2098009692989996
$ws.Range("P14").Value = 0.2098009692989996
$ws.Range("Q14").Value = 9269.034737816102
$ws.Range("R14").Value = 83421.31264034491
$ws.Range("S14").Value = 0.003344751969411691
$ws.Range("T14").Value = 0.003344751969411692

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 93.67702500000001
$ws.Range("H15").Value = 281.031075
$ws.Range("I15").Value = 0.01594250007799006
$ws.Range("J15").Value = 0.01594250007799006
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 163.0062356666667
$ws.Range("N15").Value = 489.018707
$ws.Range("O15").Value = 0.345629090707923
$ws.Range("P15").Value = 0.3456290907079231
$ws.Range("Q15").Value = 15269.93921370223
$ws.Range("R15").Value = 137429.4529233201
$ws.Range("S15").Value = 0.005510191805566697
$ws.Range("T15").Value = 0.005510191805566698

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 93.67702500000001
$ws.Range("H16").Value = 281.031075
$ws.Range("I16").Value = 0.01594250007799006
$ws.Range("J16").Value = 0.01594250007799006
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 65.39610666666668
$ws.Range("N16").Value = 196.18832
$ws.Range("O16").Value = 0.1386621609326595
$ws.Range("P16").Value = 0.1386621609326595
$ws.Range("Q16").Value = 6126.112719116002
$ws.Range("R16").Value = 55135.01447204401
$ws.Range("S16").Value = 0.002210621511483194
$ws.Range("T16").Value = 0.002210621511483194

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 93.67702500000001
$ws.Range("H17").Value = 281.031075
$ws.Range("I17").Value = 0.01594250007799006
$ws.Range("J17").Value = 0.01594250007799006
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 144.2727966666667
$ws.Range("N17").Value = 432.81839
$ws.Range("O17").Value = 0.3059077790604178
$ws.Range("P17").Value = 0.3059077790604178
$ws.Range("Q17").Value = 13515.04638016325
$ws.Range("R17").Value = 121635.4174214693
$ws.Range("S17").Value = 0.004876934791528478
$ws.Range("T17").Value = 0.004876934791528478
